$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column-D price cells to stay text (Excel would otherwise parse
# numeric-looking strings like "1.00" into the number 1 and lose formatting).
$priceCells = @("D2","D3","D5","D6","D7","D8","D9","D10","D11","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D26","D27","D29","D32","D34","D35","D36","D37","D38","D40","D43","D45","D47","D48","D49","D50","D51")
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "33.774.52"
$ws.Range("E2").Value = "  +6.76%  "
$ws.Range("D3").Value = "1.777.11"
$ws.Range("E3").Value = "  +3.83%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "225.17"
$ws.Range("E5").Value = "  +1.45%  "
$ws.Range("D6").Value = "0.559"
$ws.Range("E6").Value = "  +4.44%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").Value = "30.36"
$ws.Range("E8").Value = "  +1.42%  "
$ws.Range("D9").Value = "46.77"
$ws.Range("E9").Value = "  +3.10%  "
$ws.Range("D10").Value = "0.277"
$ws.Range("E10").Value = "  +2.59%  "
$ws.Range("D11").Value = "0.0664"
$ws.Range("E11").Value = "  +3.19%  "
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("D13").Value = "2.031.88"
$ws.Range("E13").Value = "  +3.72%  "
$ws.Range("D14").Value = "1.779.00"
$ws.Range("E14").Value = "  +3.99%  "
$ws.Range("D15").Value = "0.623"
$ws.Range("E15").Value = "  +1.77%  "
$ws.Range("D16").Value = "33.733.26"
$ws.Range("E16").Value = "  +6.81%  "
$ws.Range("D17").Value = "9.99"
$ws.Range("E17").Value = "  -2.43%  "
$ws.Range("D18").Value = "4.17"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").Value = "68.47"
$ws.Range("E19").Value = "  +1.80%  "
$ws.Range("D20").Value = "250.54"
$ws.Range("D21").Value = "0.0₃0737"
$ws.Range("E21").Value = "  +1.87%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").Value = "10.26"
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("D24").Value = "4.16"
$ws.Range("E24").Value = "  -2.40%  "
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("D26").Value = "158.27"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("D27").Value = "16.42"
$ws.Range("E27").Value = "  +2.35%  "
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("D29").Value = "6.96"
$ws.Range("E29").Value = "  +2.28%  "
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("D32").Value = "0.0514"
$ws.Range("E32").Value = "  +2.07%  "
$ws.Range("E33").Value = "  +2.93%  "
$ws.Range("D34").Value = "3.55"
$ws.Range("E34").Value = "  +4.49%  "
$ws.Range("D35").Value = "1.82"
$ws.Range("E35").Value = "  +4.88%  "
$ws.Range("D36").Value = "1.478.88"
$ws.Range("E36").Value = "  -2.57%  "
$ws.Range("D37").Value = "1.06"
$ws.Range("E37").Value = "  +2.62%  "
$ws.Range("D38").Value = "0.627"
$ws.Range("E38").Value = "  +3.00%  "
$ws.Range("E39").Value = "  +2.03%  "
$ws.Range("D40").Value = "82.97"
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("E41").Value = "  +1.94%  "
$ws.Range("E42").Value = "  -1.29%  "
$ws.Range("D43").Value = "0.884"
$ws.Range("E43").Value = "  +3.48%  "
$ws.Range("E44").Value = "  +1.80%  "
$ws.Range("D45").Value = "0.0507"
$ws.Range("E45").Value = "  +0.42%  "
$ws.Range("E46").Value = "  +3.42%  "
$ws.Range("D47").Value = "1.929.20"
$ws.Range("E47").Value = "  +4.39%  "
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  +0.25%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "5.69"
$ws.Range("E49").Value = "  +2.25%  "
$ws.Range("D50").Value = "11.85"
$ws.Range("E50").Value = "  +15.27%  "
$ws.Range("D51").Value = "50.87"
$ws.Range("E51").Value = "  -2.91%  "
